$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp shown in F1
$ws.Range("F1").Value = "Last status check on: 18.02.2022 15:15"

# New scrape results for row 2 (TankONO):
#   the previous "current price" (B2) becomes the "old price" (C2)
$ws.Range("C2").Value = 35.5
#   the freshly scraped price becomes the new "current price" (B2)
$ws.Range("B2").Value = 35.9

# The price delta is now written out as a signed text string instead of a number
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "+0.4"
$ws.Range("D2").Style = "Normal"

# The "old date" column now stores a plain text timestamp instead of a numeric date serial
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2022-02-18 15:15:11"
$ws.Range("E2").Style = "Normal"
